# Update of daily and weekly charts
# Apply the revised daily new-case counts (column C) for CasesByDate, which
# cascade into the running total (column B) and 7-day average (column D)
# through their existing formulas, then append the newest day's data (row 317).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CasesByDate")

# --- Revised daily case counts (column C) -------------------------------
$cUpdates = @{
    86  = 2704
    179 = 257
    225 = 547
    227 = 409
    246 = 720
    255 = 834
    261 = 899
    269 = 1374
    294 = 3521
    295 = 3134
    296 = 2919
    297 = 2997
    298 = 2854
    300 = 1187
    302 = 3789
    305 = 3363
    308 = 5477
    309 = 5834
    310 = 6025
    311 = 5622
    312 = 4998
    313 = 1859
    314 = 1611
    315 = 4370
    316 = 3052
}

foreach ($row in $cUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $cUpdates[$row]
}

# --- Append the new day: row 317 ----------------------------------------
$newRow = 317
$prevRow = $newRow - 1

# Carry the date formatting down from the row above (reuses the existing
# style record instead of minting a duplicate number format).
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 44174
$ws.Cells.Item($newRow, 2).Formula = "=C$newRow+B$prevRow"
$ws.Cells.Item($newRow, 3).Value = 243
$ws.Cells.Item($newRow, 4).Formula = "=AVERAGE(C$($newRow-6):C$newRow)"

$excel.Calculate()
